# Update crypto price/volume data per commit "Updated cryptos list on Wed Sep 11 19:21:08 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.067.70"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").Value = "2.315.81"
$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'533.17"
$ws.Range("E5").Value = "  +2.56%  "

$ws.Range("D6").Value = "'132.10"
$ws.Range("E6").Value = "  -2.94%  "

$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("D8").Value = "'0.537"
$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("D9").Value = "2.337.09"
$ws.Range("E9").Value = "  -1.03%  "

$ws.Range("D10").Value = "'0.102"
$ws.Range("E10").Value = "  -1.00%  "

$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("D12").Value = "'5.29"
$ws.Range("E12").Value = "  -2.40%  "

$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").Value = "2.750.03"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").Value = "'23.46"
$ws.Range("E15").Value = "  -2.97%  "

$ws.Range("D16").Value = "57.118.05"
$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("E17").Value = "  -2.05%  "

$ws.Range("D18").Value = "2.333.31"
$ws.Range("E18").Value = "  -0.89%  "

$ws.Range("D19").Value = "'338.14"
$ws.Range("E19").Value = "  +2.85%  "

$ws.Range("D20").Value = "'10.46"
$ws.Range("E20").Value = "  -1.27%  "

$ws.Range("E21").Value = "  +2.96%  "

$ws.Range("D22").Value = "'4.15"
$ws.Range("E22").Value = "  -2.01%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").Value = "'61.62"
$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("D25").Value = "'8.75"
$ws.Range("E25").Value = "  +5.97%  "

$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("D28").Value = "'1.33"
$ws.Range("E28").Value = "  +0.37%  "

$ws.Range("D29").Value = "'170.66"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("E30").Value = "  +1.35%  "

$ws.Range("D31").Value = "0.0₃0721"
$ws.Range("E31").Value = "  -2.91%  "

$ws.Range("D32").Value = "'6.09"
$ws.Range("E32").Value = "  -2.73%  "

$ws.Range("D33").Value = "'18.51"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("E35").Value = "  -0.35%  "

$ws.Range("E36").Value = "  -2.70%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.00"
$ws.Range("E37").Value = "  -0.98%  "

$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "'0.905"
$ws.Range("E38").Value = "  -1.77%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.58"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'39.11"
$ws.Range("E40").Value = "  +1.51%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'5.74"
$ws.Range("E41").Value = "  +9.15%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'148.71"
$ws.Range("E42").Value = "  -1.64%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.377"
$ws.Range("E43").Value = "  -1.43%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'3.59"
$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'280.17"
$ws.Range("E45").Value = "  -0.55%  "

$ws.Range("D46").Value = "'0.0928"
$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("D47").Value = "'0.0501"
$ws.Range("E47").Value = "  -1.36%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'18.62"
$ws.Range("E48").Value = "  +2.22%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.555"
$ws.Range("E49").Value = "  -1.26%  "

$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D50").Value = "'0.382"
$ws.Range("E50").Value = "  +0.21%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0216"
$ws.Range("E51").Value = "  -2.09%  "
